$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (A205, which carries
# the bold/bordered/centered style used throughout column A) down onto the
# new rows 206:217 before writing values into them.
$ws.Range("A205").Copy()
$ws.Range("A206:A217").PasteSpecial(-4122)

$data = @(
    @(204, 0.7700000000000001),
    @(205, 0.5300000000000001),
    @(206, 0.7220000000000001),
    @(207, 0.3700000000000002),
    @(208, 0.4100000000000001),
    @(209, 0.2557142857142859),
    @(210, 0.3700000000000002),
    @(211, 0.7220000000000001),
    @(212, 0.5300000000000001),
    @(213, 0.5300000000000001),
    @(214, 0.5300000000000001),
    @(215, 0.5300000000000001)
)

$row = 206
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
